$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("1-9")
$ws2 = $wb.Worksheets.Item("10")

# --- Fill in the new survey respondent rows (row 14 = #10, row 15 = #11) ---
$ws1.Range("B14").Value = "Nguyễn Minh Thảo"
$ws1.Range("C14").Value = "lightning_spkt@yahoo.com.vn"
$ws1.Range("D14").Value = 4
$ws1.Range("E14").Value = 3
$ws1.Range("F14").Value = 4
$ws1.Range("G14").Value = 4
$ws1.Range("H14").Value = 3
$ws1.Range("I14").Value = 4
$ws1.Range("J14").Value = 4
$ws1.Range("K14").Value = 3
$ws1.Range("L14").Value = 4

$ws1.Range("B15").Value = "Nguyễn Chí Hiếu"
$ws1.Range("C15").Value = "asakura255@gmail.com"
$ws1.Range("D15").Value = 2
$ws1.Range("E15").Value = 2
$ws1.Range("F15").Value = 2
$ws1.Range("G15").Value = 2
$ws1.Range("H15").Value = 3
$ws1.Range("I15").Value = 3
$ws1.Range("J15").Value = 2
$ws1.Range("K15").Value = 2
$ws1.Range("L15").Value = 2

# Turn the two new e-mail addresses into live mailto: hyperlinks, same as
# the other respondents' e-mail cells already have.
$ws1.Hyperlinks.Add($ws1.Range("C14"), "mailto:lightning_spkt@yahoo.com.vn")
$ws1.Hyperlinks.Add($ws1.Range("C15"), "mailto:asakura255@gmail.com")

# --- Update the saved view/selection state ---
# Sheet "10" is no longer the active tab; stop scrolled near the bottom and
# keep its previous selection untouched.
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 36
[void]$ws2.Range("A43:N46").Select()

# Sheet "1-9" becomes the active tab, scrolled so row 7 is at the top, with
# C20 as the selected cell.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("C20").Select()
